$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H2, H4 and H6 hold the "vOutData" timestamp for CT 01 / CT 03 / CT 05.
# They all get updated to the same new text value "12/05/2020".
# Using Formula + Copy/PasteSpecial(values) writes the text as a literal
# string (instead of a date serial, which is what a plain .Value=
# assignment would infer from a dd/mm/yyyy-looking string) while leaving
# the cell's existing style untouched.
$cells = @("H2", "H4", "H6")
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    $rng.Formula = '="12/05/2020"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
